# Week 15 simulations - append new week's per-play data to the YDS and ST
# shared-string "play lists", and roll the season-to-date totals forward on
# the OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: per-play yardage lists (space separated numbers kept as text)
# ---------------------------------------------------------------------
$yds = $wb.Worksheets.Item("YDS")

$yds.Range("B2").Value = "-1 1 6 6 3 3 5 3 3 6 7 1 3 9 4 1 4 8 6 5 4 3 8 3 1 -4 5 0 3 4 1 8 -1 3 5 2 5 -4 1 3 2 4 2 10 4 2 6 1 6 3 7 2 6 2 12 4 2 12 4 2 4 4 3 0 7 3 6 4 4 -2 15 9 6 5 3 6 3 1 -3 3 25 1 -1 5 4 5 6 2 -1 3 3 14 11 3 0 2 1 3 8 1 5 7 17 -1 5 4 1 1 57 0 -1 -4 -1 16 6 0 3 3 0 5 8 2 6 36 1 -6 28 1 5 1 5 1 0 7 6 8 3 9 3 3 9 15 2 2 13 2 1 5 3 2 0 2 0 -1 2 3 7 4 12 6 5 8 2 1 11 4 1 3 8 6 0 10 -2 7 0 9 4 1 4 2 4 8 16 0 0 4 4 9 15 2 3 -1 7 7 -2 8 5 3 1 5 5 0 7 8 2 2 4 9 6 4 9 3 11 5 -1 3 7 5 8 1 2 0 2 0 5 4 7 2 1 6 2 2 5 4 3 4 2 2 5 2 0 2 2 4 1 5 3 2 1 0 11 2 5 4 18 11 7 6 5 7 4 3 6 2 4 1 -1 0 1 2 7 5 0 8 4 1 8 2 5 2 4 -4 3 2 5 -1 2 7 7 3 3 15 1 2 3 6 11 0 11 5 8 11 1 3 14 0 5 3 6 6 2 1 9"

$yds.Range("B3").Value = "4 6 7 14 4 31 10 19 7 8 5 6 2 7 18 19 4 5 6 9 13 1 10 50 22 14 3 9 5 22 11 11 0 15 8 0 42 11 -4 2 12 2 47 14 1 0 10 3 4 15 19 6 1 9 6 7 12 25 17 5 21 7 16 4 8 12 6 6 23 7 19 10 25 26 11 1 28 2 11 5 15 17 -4 34 14 0 12 13 24 -1 3 15 24 5 22 5 7 0 8 8 2 2 59 20 20 15 13 4 0 6 32 1 12 5 11 9 4 10 12 14 13 41 8 15 3 4 12 5 0 12 1 10 3 15 11 7 13 6 19 5 10 2 20 4 10 26 4 25 17 9 7 7 12 4 1 4 15 -2 7 2 19 33 8 5 6 9 12 10 4 7 10 15 5 6 21 35 0 11 -1 14 16 3 8 2 15 5 3 41 9 4 6 12 6 34 3 4 14 8 10 6 4 23 24 4 7 8 5 50 7 37 2 4 13 3 5 39 8 17 12 15 8 25 13 12 17 10 9 13 26 4 18 75 8 15 0 2 7 7 5 1 54 27 43 7 10 6 7 8 28 6 6 5 2 9 2 8 6 7 9 6 6 32 0 7 12 25 12 20 6 13 6 18 38 12 16 8 23 16 12 6 6 8 11 6 4 7 3"

$yds.Range("C2").Value = "1 11 5 15 3 3 6 2 8 2 10 5 3 9 5 9 3 11 1 1 6 6 14 2 -1 3 2 1 1 1 2 7 4 0 2 6 3 1 4 2 8 5 2 6 3 3 8 3 1 12 2 17 3 2 0 26 4 -3 0 9 2 2 0 1 2 2 16 8 4 0 7 9 5 1 1 3 -1 10 0 2 4 5 4 8 -2 2 3 11 6 1 3 4 1 0 7 10 -1 9 4 0 6 5 7 11 1 -1 4 9 8 2 3 -2 10 4 4 1 9 2 0 0 1 25 5 4 1 3 7 3 4 2 13 4 2 2 11 3 6 1 2 14 3 16 3 8 5 4 5 2 2 6 1 4 38 2 6 4 5 3 0 2 6 27 12 5 2 8 7 13 4 2 3 1 11 -1 4 5 4 2 2 2 7 6 9 0 0 9 1 3 7 7 2 6 0 0 6 4 7 7 7 2 4 10 2 -4 -1 -1 2 3 2 2 7 10 6 0 5 4 5 1 4 7 1 12 5 2 4 -1 6 -1 3 3 1 2 5 0 -1 4 0 7 4 3 5 13 2 4 0 7 2 0 4 3 7 12 5 6 6 3 1 2 3 7 5 7 2 0 3 5 3 5 6 0 5 8 2 1 5 1 20 1 1 3 8 3 8 7 0 20 4 5 2 2 3 11 17"

$yds.Range("C3").Value = "14 3 17 2 3 3 8 5 6 4 10 8 55 8 46 9 5 6 4 5 11 4 10 8 20 1 6 15 7 2 5 3 5 7 7 5 24 16 7 8 7 6 4 16 8 16 10 11 5 7 19 1 0 16 8 0 16 12 9 3 39 12 12 12 3 9 2 45 6 8 10 4 17 4 6 12 11 11 5 0 4 11 10 7 9 30 3 4 2 7 9 6 7 2 10 4 7 5 70 28 14 5 8 3 10 2 7 19 9 8 8 5 21 15 11 20 3 5 3 12 8 11 20 21 8 10 5 15 7 6 12 40 0 5 5 14 4 17 28 9 12 5 7 10 3 11 2 19 10 3 12 13 6 15 4 8 58 6 9 8 3 6 8 11 20 9 9 4 9 23 29 9 6 22 2 7 10 4 2 1 7 2 7 12 2 25 6 14 10 13 22 9 13 12 -6 6 28 1 9 9 7 5 7 18 8 6 -4 4 8 10 8 9 5 43 3 18 56 15 10 5 7 15 6 7 9 21 5 19 6 7 7 3 23 6 19 26 4 5 79 5 11 22 3 19 6 1 8 9 54 13 7 10 7 7 7 11 14 19 8 3 6 46 6 54 6 2 5 6 13 5 11 4 3 22 5"

# ---------------------------------------------------------------------
# OFF sheet: season-to-date offensive totals, Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$off = $wb.Worksheets.Item("OFF")

$off.Range("C2").Value = 167
$off.Range("D2").Value = 8
$off.Range("F2").Value = 46
$off.Range("G2").Value = 49
$off.Range("L2").Value = 217
$off.Range("M2").Value = 149
$off.Range("O2").Value = 17
$off.Range("P2").Value = 10
$off.Range("Q2").Value = 391

$off.Range("B3").Value = 10
$off.Range("C3").Value = 130
$off.Range("D3").Value = 6
$off.Range("E3").Value = 19
$off.Range("F3").Value = 88
$off.Range("G3").Value = 35
$off.Range("H3").Value = 24
$off.Range("I3").Value = 47
$off.Range("J3").Value = 47
$off.Range("N3").Value = 14

# ---------------------------------------------------------------------
# DEF sheet: season-to-date defensive totals, Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$def = $wb.Worksheets.Item("DEF")

$def.Range("C2").Value = 143
$def.Range("D2").Value = 8
$def.Range("F2").Value = 41
$def.Range("G2").Value = 40
$def.Range("I2").Value = 4
$def.Range("J2").Value = 32
$def.Range("L2").Value = 224
$def.Range("M2").Value = 136
$def.Range("O2").Value = 21
$def.Range("Q2").Value = 363

$def.Range("B3").Value = 10
$def.Range("C3").Value = 138
$def.Range("E3").Value = 23
$def.Range("F3").Value = 87
$def.Range("G3").Value = 36
$def.Range("H3").Value = 11
$def.Range("I3").Value = 40
$def.Range("J3").Value = 48
$def.Range("N3").Value = 14

# ---------------------------------------------------------------------
# ST sheet: special teams totals + per-kick distance lists
# ---------------------------------------------------------------------
$st = $wb.Worksheets.Item("ST")

$st.Range("B2").Value = 67
$st.Range("D2").Value = 40
$st.Range("F2").Value = 704
$st.Range("G2").Value = 687
$st.Range("J2").Value = 253
$st.Range("K2").Value = 230

$st.Range("B3").Value = 30

$st.Range("D3").Value = "38 47 59 32 46 56 38 51 50 63 36 57 45 36 31 36 50 82 49 53 37 58 48 39 55 60 50 41 49 40 58 23 41 53 61 34 57 45 47 22"
$st.Range("B4").Value = "57 63 65 64 66 73 70 69 69 64 68 66 62 60 67 65 66 66 67 65 70 68 43 57 63 63 64 44 62 66 60 64 61 65 63 66 62"
$st.Range("D4").Value = "0 0 13 0 0 10 0 0 9 14 0 0 14 2 0 11 27 0 0 0 0 0 14 3 0 0 0 0 0 0 8 0 0 13 0 0 34 0 97 0"
$st.Range("B5").Value = "33 21 68 16 20 38 27 40 44 22 28 26 26 27 30 27 30 24 25 23 20 32 12 12 24 32 23 15 23 20 28 40 0 42 22 14 18"
$st.Range("D5").Value = "0 17 2 0 0 0 0 6 9 0 0 1 0 0 8 0 16 0 0 0 8 17 0 0 0 0 15 9 0 0 4 0 16 0 0 0 0 0 0 0 0 13 0 15"
$st.Range("B6").Value = "17 41 33 21 21 16 10 16 23 13 26 19 3 32 14 10 14 19"

# ---------------------------------------------------------------------
# TURNS sheet: season-to-date turnovers, Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$turns = $wb.Worksheets.Item("TURNS")

$turns.Range("C2").Value = 8
$turns.Range("D2").Value = 6
$turns.Range("E2").Value = 9

$turns.Range("D3").Value = 5
$turns.Range("E3").Value = 8

# ---------------------------------------------------------------------
# PEN sheet: season-to-date penalties
# ---------------------------------------------------------------------
$pen = $wb.Worksheets.Item("PEN")

$pen.Range("B3").Value = 12
$pen.Range("D4").Value = 9
